# Atualização de bases das ligas, do dia: 15-06-2024 às 21:10
# Cyclic update of match rows 117-120 (ids 115-118) on "Uruguay Primera División" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 117 (id 115)
$ws.Range("B117").Value = 7013885
$ws.Range("E117").Value = "La Luz"
$ws.Range("F117").Value = "Atletico Fenix Montevideo"
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 2
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = "A"
$ws.Range("L117").Value = 3
$ws.Range("M117").Value = 3
$ws.Range("N117").Value = 2.4
$ws.Range("O117").Value = 2.9
$ws.Range("P117").Value = 2.75
$ws.Range("Q117").Value = 2.6
$ws.Range("R117").Value = 0
$ws.Range("S117").Value = 2.025
$ws.Range("T117").Value = 1.825
$ws.Range("U117").Value = 2
$ws.Range("V117").Value = 2.025
$ws.Range("W117").Value = 1.825
$ws.Range("X117").Value = -1
$ws.Range("Y117").Value = -1
$ws.Range("Z117").Value = 1.6
$ws.Range("AA117").Value = -1
$ws.Range("AB117").Value = 0.825
$ws.Range("AC117").Value = 0
$ws.Range("AD117").Value = 0

# Row 118 (id 116)
$ws.Range("B118").Value = 7013702
$ws.Range("E118").Value = "Defensor Sporting"
$ws.Range("F118").Value = "Danubio"
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 2
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 1
$ws.Range("K118").Value = "A"
$ws.Range("L118").Value = 1.8
$ws.Range("M118").Value = 3.6
$ws.Range("N118").Value = 4.2
$ws.Range("O118").Value = 1.8
$ws.Range("P118").Value = 3.6
$ws.Range("Q118").Value = 4.2
$ws.Range("R118").Value = -0.75
$ws.Range("S118").Value = 2.05
$ws.Range("T118").Value = 1.8
$ws.Range("U118").Value = 2.25
$ws.Range("V118").Value = 1.85
$ws.Range("W118").Value = 2
$ws.Range("X118").Value = -1
$ws.Range("Y118").Value = -1
$ws.Range("Z118").Value = 3.2
$ws.Range("AA118").Value = -1
$ws.Range("AB118").Value = 0.8
$ws.Range("AC118").Value = -0.5
$ws.Range("AD118").Value = 0.5

# Row 119 (id 117)
$ws.Range("B119").Value = 7013409
$ws.Range("E119").Value = "Nacional De Football"
$ws.Range("F119").Value = "Torque"
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 1
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 1
$ws.Range("K119").Value = "D"
$ws.Range("L119").Value = 1.666
$ws.Range("M119").Value = 3.9
$ws.Range("N119").Value = 4.5
$ws.Range("O119").Value = 1.615
$ws.Range("P119").Value = 4
$ws.Range("Q119").Value = 4.75
$ws.Range("R119").Value = -0.75
$ws.Range("S119").Value = 1.8
$ws.Range("T119").Value = 2.05
$ws.Range("U119").Value = 2.75
$ws.Range("V119").Value = 1.95
$ws.Range("W119").Value = 1.9
$ws.Range("X119").Value = -1
$ws.Range("Y119").Value = 3
$ws.Range("Z119").Value = -1
$ws.Range("AA119").Value = -1
$ws.Range("AB119").Value = 1.05
$ws.Range("AC119").Value = -1
$ws.Range("AD119").Value = 0.8999999999999999

# Row 120 (id 118)
$ws.Range("B120").Value = 7013886
$ws.Range("E120").Value = "Racing Club de Montevideo"
$ws.Range("F120").Value = "Cerro"
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 1
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = "A"
$ws.Range("L120").Value = 2.25
$ws.Range("M120").Value = 3.1
$ws.Range("N120").Value = 3.25
$ws.Range("O120").Value = 2.25
$ws.Range("P120").Value = 2.875
$ws.Range("Q120").Value = 3.5
$ws.Range("R120").Value = -0.25
$ws.Range("S120").Value = 1.95
$ws.Range("T120").Value = 1.9
$ws.Range("U120").Value = 2
$ws.Range("V120").Value = 1.925
$ws.Range("W120").Value = 1.925
$ws.Range("X120").Value = -1
$ws.Range("Y120").Value = -1
$ws.Range("Z120").Value = 2.5
$ws.Range("AA120").Value = -1
$ws.Range("AB120").Value = 0.8999999999999999
$ws.Range("AC120").Value = -1
$ws.Range("AD120").Value = 0.925
